# The "Search" worksheet holds label/value pairs describing the search
# screen. A new "search_full_text" criterion row needs to be inserted
# right after the existing search rows (date / by / nb_results / criteria)
# and before the streetName/city rows, pushing those down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")

# Insert a new blank row at row 4 - this shifts the existing row 4
# (address_streetName / streetName) down to row 5, and the existing
# row 5 (address_city / city) down to row 6.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the search_full_text label/value.
$ws.Range("A4").Value = "`${msg.getProperty('search_full_text')}"
$ws.Range("B4").Value = "`${search_full_text}"
